$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1252.478759765625
$ws.Range("C2").Value = 0.9396
$ws.Range("D2").Value = 0.9196000099182129
$ws.Range("E2").Value = 1.358000040054321
$ws.Range("F2").Value = 0.7696999907493591
$ws.Range("H2").Value = 0.7351
$ws.Range("B3").Value = 1143.126831054688
$ws.Range("C3").Value = 0.9145
$ws.Range("D3").Value = 0.9062
$ws.Range("E3").Value = 1.104900002479553
$ws.Range("F3").Value = 0.7699999809265137
$ws.Range("H3").Value = 0.6163999999999999
$ws.Range("B4").Value = 787.7871704101562
$ws.Range("C4").Value = 0.9203
$ws.Range("D4").Value = 0.9111
$ws.Range("E4").Value = 1.07860004901886
$ws.Range("F4").Value = 0.86080002784729
$ws.Range("H4").Value = 0.66
$ws.Range("B5").Value = 866.4793090820312
$ws.Range("C5").Value = 0.914
$ws.Range("D5").Value = 0.9165
$ws.Range("E5").Value = 0.9872000217437744
$ws.Range("F5").Value = 0.6866999864578247
$ws.Range("H5").Value = 0.7077
$ws.Range("B6").Value = 1152.57763671875
$ws.Range("C6").Value = 0.9118000000000001
$ws.Range("D6").Value = 0.9131
$ws.Range("E6").Value = 0.960099995136261
$ws.Range("F6").Value = 0.809499979019165
$ws.Range("H6").Value = 0.6777
$ws.Range("B7").Value = 906.7448120117188
$ws.Range("C7").Value = 0.9131
$ws.Range("D7").Value = 0.9132000207901001
$ws.Range("E7").Value = 0.9473000168800354
$ws.Range("F7").Value = 0.8410000205039978
$ws.Range("H7").Value = 0.6788999999999999
$ws.Range("B8").Value = 1014.243530273438
$ws.Range("C8").Value = 0.9088000000000001
$ws.Range("D8").Value = 0.9088000000000001
$ws.Range("E8").Value = 0.9652000069618225
$ws.Range("F8").Value = 0.8661999702453613
$ws.Range("H8").Value = 0.6399
$ws.Range("B9").Value = 7123.4375
$ws.Range("C9").Value = 0.918
$ws.Range("D9").Value = 0.9121
$ws.Range("E9").Value = 1.358000040054321
$ws.Range("F9").Value = 0.6866999864578247
$ws.Range("H9").Value = 4.715699999999999
